$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new shared string used by row 4 (sending cluster renamed)
$ws.Range("A4").Value = "Inflammatory-Mac"

# Row 2 updated numeric values
$ws.Range("G2").Value = 0.07908133333333334
$ws.Range("H2").Value = 0.237244
$ws.Range("I2").Value = 0.1380838080781507
$ws.Range("J2").Value = 0.1380838080781507
$ws.Range("Q2").Value = 0.004934279793333334
$ws.Range("R2").Value = 0.04440851814
$ws.Range("S2").Value = 0.1380838080781507
$ws.Range("T2").Value = 0.1380838080781507

# Row 3 updated numeric values
$ws.Range("I3").Value = 0.1348139473702591
$ws.Range("J3").Value = 0.134813947370259
$ws.Range("S3").Value = 0.1348139473702591
$ws.Range("T3").Value = 0.134813947370259

# Row 4 updated numeric values
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03917266666666667
$ws.Range("H4").Value = 0.117518
$ws.Range("I4").Value = 0.06839933974190335
$ws.Range("J4").Value = 0.06839933974190333
$ws.Range("Q4").Value = 0.002444178536666667
$ws.Range("R4").Value = 0.02199760683
$ws.Range("S4").Value = 0.06839933974190335
$ws.Range("T4").Value = 0.06839933974190333

# New row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Ccl21b"
$ws.Range("C5").Value = "Ackr2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3772426666666667
$ws.Range("H5").Value = 1.131728
$ws.Range("I5").Value = 0.658702904809687
$ws.Range("J5").Value = 0.6587029048096869
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.062395
$ws.Range("N5").Value = 0.187185
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.02353805618666667
$ws.Range("R5").Value = 0.21184250568
$ws.Range("S5").Value = 0.658702904809687
$ws.Range("T5").Value = 0.6587029048096869
